# Append three new history entries after the last existing one
# ("... Dorfuchs etc."):
#   "Informationen zu allen Unterseiten + Videos"
#   "14.03. "
#   "Icons von Leonie und Phillip auf Hauptseite"
# The trailing hidden "_GoBack" bookmark must end up wrapping the very end
# of the new last paragraph (after its run), exactly as it did around the
# old last paragraph before this edit.

$d = $word.ActiveDocument

# Remember which paragraph currently ends in "etc." (the last "real"
# paragraph of the document) before we touch anything.
$probe = $d.Content
$probe.Find.Execute("etc.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$etcParaIndex = $probe.Paragraphs.Item(1).Index

# Split that paragraph's tail into four paragraphs via Find/Replace using
# "^p" paragraph marks. Word keeps the hidden bookmark glued to the old
# end-of-match position, so after this call it has moved into a brand new,
# otherwise-empty trailing paragraph instead of staying inside the
# "... etc." paragraph - exactly the 3rd paragraph after $etcParaIndex.
$r = $d.Content
$r.Find.Execute("etc.", $true, $false, $false, $false, $false, $true, 1, $false, `
    "etc.^pInformationen zu allen Unterseiten + Videos^p14.03. ^p", 2) | Out-Null

# Insert the final line of text ahead of the (now isolated) bookmark, in
# the same paragraph, so the bookmark ends up trailing the run instead of
# leading it - matching how it originally trailed the "... etc." run.
$bookmarkParaIndex = $etcParaIndex + 3
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)
$bookmarkPara.Range.InsertBefore("Icons von Leonie und Phillip auf Hauptseite")
